# Rebuild the "109_2" nominations summary sheet with the new, more granular
# row layout (per-category breakdown + renamed/re-ordered totals) and drop
# the old trailing "Total Returned to the White House" row (old row 44).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Labels"
$ws.Range("B1").Value = "Values"

$ws.Range("A2").Value = "Congress"
$ws.Range("B2").Value = 109

$ws.Range("A3").Value = "Session"
$ws.Range("B3").Value = 2

$ws.Range("A4").Value = "Start Date"
$ws.Range("B4").Value = 38720

$ws.Range("A5").Value = "End Date"
$ws.Range("B5").Value = 39082

$ws.Range("A6").Value = "Civilian "

$ws.Range("A7").Value = "     Civilian, New nominations"
$ws.Range("B7").Value = 470

$ws.Range("A8").Value = "     Civilian, Carryover nominations"
$ws.Range("B8").Value = 148

$ws.Range("A9").Value = "     Civilian, Confirmed "
$ws.Range("B9").Value = 415

$ws.Range("A10").Value = "     Civilian, Withdrawn "
$ws.Range("B10").Value = 21

$ws.Range("A11").Value = "     Civilian, Returned to White House "
$ws.Range("B11").Value = 182

$ws.Range("A12").Value = "Other Civilian "

$ws.Range("A13").Value = "     Other Civilian, New nominations"
$ws.Range("B13").Value = 2486

$ws.Range("A14").Value = "     Other Civilian, Carryover nominations"
$ws.Range("B14").Value = 780

$ws.Range("A15").Value = "     Other Civilian, Confirmed "
$ws.Range("B15").Value = 3263
$ws.Range("B15").NumberFormat = "#,##0"
$ws.Range("B15").HorizontalAlignment = -4152

$ws.Range("A16").Value = "     Other Civilian, Withdrawn "
$ws.Range("B16").Value = 1

$ws.Range("A17").Value = "     Other Civilian, Returned to White House "
$ws.Range("B17").Value = 2

$ws.Range("A18").Value = "Air Force "

$ws.Range("A19").Value = "     Air Force, New nominations"
$ws.Range("B19").Value = 7730

$ws.Range("A20").Value = "     Air Force, Carryover nominations"
$ws.Range("B20").Value = 100

$ws.Range("A21").Value = "     Air Force, Confirmed "
$ws.Range("B21").Value = 7829
$ws.Range("B21").NumberFormat = "#,##0"
$ws.Range("B21").HorizontalAlignment = -4152

$ws.Range("A22").Value = "     Air Force, Returned to White House "
$ws.Range("B22").Value = 1

$ws.Range("A23").Value = "Army "

$ws.Range("A24").Value = "     Army, New nominations"
$ws.Range("B24").Value = 9177

$ws.Range("A25").Value = "     Army, Carryover nominations"
$ws.Range("B25").Value = 608

$ws.Range("A26").Value = "     Army, Confirmed "
$ws.Range("B26").Value = 9772
$ws.Range("B26").NumberFormat = "#,##0"
$ws.Range("B26").HorizontalAlignment = -4152

$ws.Range("A27").Value = "     Army, Returned to White House "
$ws.Range("B27").Value = 13

$ws.Range("A28").Value = "Navy "

$ws.Range("A29").Value = "     Navy, New nominations"
$ws.Range("B29").Value = 7015

$ws.Range("A30").Value = "     Navy, Carryover nominations"
$ws.Range("B30").Value = 21

$ws.Range("A31").Value = "     Navy, Confirmed "
$ws.Range("B31").Value = 7035
$ws.Range("B31").NumberFormat = "#,##0"
$ws.Range("B31").HorizontalAlignment = -4152

$ws.Range("A32").Value = "     Navy, Returned to White House "
$ws.Range("B32").Value = 1

$ws.Range("A33").Value = "Marine Corps "

$ws.Range("A34").Value = "     Marine Corps, New nominations"
$ws.Range("B34").Value = 1291

$ws.Range("A35").Value = "     Marine Corps, Carryover nominations"
$ws.Range("B35").Value = 2

$ws.Range("A36").Value = "     Marine Corps, Confirmed "
$ws.Range("B36").Value = 1289
$ws.Range("B36").NumberFormat = "#,##0"
$ws.Range("B36").HorizontalAlignment = -4152

$ws.Range("A37").Value = "     Marine Corps, Returned to White House "
$ws.Range("B37").Value = 4

$ws.Range("A38").Value = "Total new nominations"
$ws.Range("B38").Value = 28169
$ws.Range("B38").NumberFormat = "#,##0"
$ws.Range("B38").HorizontalAlignment = -4152

$ws.Range("A39").Value = "Total carryover nominations"
$ws.Range("B39").Value = 1659
$ws.Range("B39").NumberFormat = "#,##0"
$ws.Range("B39").HorizontalAlignment = -4152

$ws.Range("A40").Value = "Total confirmed"
$ws.Range("B40").Value = 29603
$ws.Range("B40").NumberFormat = "#,##0"
$ws.Range("B40").HorizontalAlignment = -4152

$ws.Range("A41").Value = "Total unconfirmed"
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B41").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B41").Value = 0

$ws.Range("A42").Value = "Total withdrawn"
$ws.Range("B42").Value = 22

$ws.Range("A43").Value = "Total returned"
$ws.Range("B43").Value = 203

# The old sheet had one extra trailing row ("Total Returned to the White
# House") that no longer exists in the new layout - remove it so the used
# range shrinks back down to A1:B43.
$ws.Rows.Item(44).Delete()